# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue {
    param($Worksheet, $CellRef, $Text)
    # Force the value to be stored as text even when it looks numeric
    # (e.g. "556.39"), then restore the default "General"/Normal style
    # so no stray per-cell formatting is introduced.
    $range = $Worksheet.Range($CellRef)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "67.647.64"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "3.519.82"
$ws.Range("E3").Value = "  -3.18%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextCellValue $ws "D5" "203.11"
$ws.Range("E5").Value = "  +2.75%  "
Set-TextCellValue $ws "D6" "556.39"
$ws.Range("E6").Value = "  -3.68%  "
$ws.Range("D7").Value = "3.508.28"
$ws.Range("E7").Value = "  -3.30%  "
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCellValue $ws "D10" "64.32"
$ws.Range("E10").Value = "  +13.66%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCellValue $ws "D11" "0.658"
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("E12").Value = "  -6.24%  "
$ws.Range("E13").Value = "  -6.75%  "
Set-TextCellValue $ws "D14" "9.91"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").Value = "4.077.29"
$ws.Range("D16").Value = "3.516.85"
$ws.Range("E16").Value = "  -3.49%  "
$ws.Range("E17").Value = "  -1.84%  "
Set-TextCellValue $ws "D18" "18.55"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "67.392.32"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("E20").Value = "  -5.66%  "
$ws.Range("E21").Value = "  -5.35%  "
Set-TextCellValue $ws "D22" "393.65"
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("E23").Value = "  -6.77%  "
$ws.Range("E24").Value = "  -5.53%  "
Set-TextCellValue $ws "D25" "83.37"
$ws.Range("E25").Value = "  -2.97%  "
Set-TextCellValue $ws "D26" "3.93"
$ws.Range("E26").Value = "  +1.26%  "
Set-TextCellValue $ws "D27" "2.84"
$ws.Range("E27").Value = "  -4.07%  "
Set-TextCellValue $ws "D28" "12.25"
$ws.Range("E28").Value = "  -3.09%  "
Set-TextCellValue $ws "D29" "8.88"
$ws.Range("E29").Value = "  -3.36%  "
Set-TextCellValue $ws "D30" "715.20"
$ws.Range("E30").Value = "  +4.08%  "
Set-TextCellValue $ws "D31" "31.13"
$ws.Range("E31").Value = "  -1.93%  "
Set-TextCellValue $ws "D32" "7.09"
$ws.Range("E32").Value = "  -13.40%  "
Set-TextCellValue $ws "D33" "11.78"
$ws.Range("E33").Value = "  -3.78%  "
Set-TextCellValue $ws "D34" "64.05"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E35").Value = "  -5.35%  "
Set-TextCellValue $ws "D36" "38.71"
$ws.Range("E36").Value = "  -9.50%  "
$ws.Range("E37").Value = "  -0.06%  "
Set-TextCellValue $ws "D38" "0.399"
$ws.Range("E38").Value = "  -6.08%  "
$ws.Range("E39").Value = "  -4.07%  "
Set-TextCellValue $ws "D40" "3.01"
$ws.Range("E40").Value = "  -4.37%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "3.060.12"
$ws.Range("E42").Value = "  -4.65%  "
$ws.Range("D43").Value = "0.0₃0689"
$ws.Range("E43").Value = "  -12.49%  "
$ws.Range("E44").Value = "  -10.46%  "
$ws.Range("E45").Value = "  +5.76%  "
Set-TextCellValue $ws "D46" "2.74"
$ws.Range("E46").Value = "  -9.86%  "
Set-TextCellValue $ws "D47" "0.0408"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("E48").Value = "  -3.29%  "
Set-TextCellValue $ws "D49" "138.56"
$ws.Range("E49").Value = "  -2.03%  "
Set-TextCellValue $ws "D50" "8.28"
$ws.Range("E50").Value = "  -7.26%  "
Set-TextCellValue $ws "D51" "2.88"
$ws.Range("E51").Value = "  -7.75%  "
